$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.425.64'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '1.824.70'
$ws.Range('E3').Value = '  -0.38%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '314.54'
$ws.Range('E5').Value = '  -0.99%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').Value = '0.5099'
$ws.Range('E7').Value = '  -4.25%  '
$ws.Range('D8').Value = '0.3927'
$ws.Range('E8').Value = '  -3.08%  '
$ws.Range('D9').Value = '0.07723'
$ws.Range('E9').Value = '  +1.57%  '
$ws.Range('D10').Value = '41.88'
$ws.Range('E10').Value = '  +0.12%  '
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('D12').Value = '21.01'
$ws.Range('E12').Value = '  +0.36%  '
$ws.Range('D13').Value = '6.259'
$ws.Range('E13').Value = '  -1.68%  '
$ws.Range('D14').Value = '1.002'
$ws.Range('E14').Value = '  +0.11%  '
$ws.Range('D15').Value = '7.509'
$ws.Range('E15').Value = '  -0.51%  '
$ws.Range('D16').Value = '1.827.98'
$ws.Range('E16').Value = '  +0.25%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = '0.00001146'
$ws.Range('E17').Value = '  +6.73%  '
$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D18').Value = '92.92'
$ws.Range('E18').Value = '  +3.99%  '
$ws.Range('D19').Value = '0.06642'
$ws.Range('E19').Value = '  +0.55%  '
$ws.Range('D20').Value = '17.73'
$ws.Range('E20').Value = '  +0.66%  '
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('D22').Value = '6.095'
$ws.Range('E22').Value = '  +0.32%  '
$ws.Range('D23').Value = '28.467.45'
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('D25').Value = '2.256'
$ws.Range('E25').Value = '  +4.52%  '
$ws.Range('D26').Value = '21.09'
$ws.Range('E26').Value = '  +2.56%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = '155.80'
$ws.Range('E27').Value = '  -0.78%  '
$ws.Range('B28').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C28').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D28').Value = '2.031.84'
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('D29').Value = '2.393'
$ws.Range('E29').Value = '  -3.81%  '
$ws.Range('D30').Value = '124.74'
$ws.Range('E30').Value = '  +0.85%  '
$ws.Range('D31').Value = '0.1099'
$ws.Range('E31').Value = '  +0.38%  '
$ws.Range('D32').Value = '1.107'
$ws.Range('E32').Value = '  -1.62%  '
$ws.Range('D33').Value = '5.671'
$ws.Range('E33').Value = '  -0.29%  '
$ws.Range('D34').Value = '3.654'
$ws.Range('E34').Value = '  -0.15%  '
$ws.Range('D35').Value = '0.07032'
$ws.Range('E35').Value = '  -2.22%  '
$ws.Range('D36').Value = '0.2214'
$ws.Range('E36').Value = '  -2.08%  '
$ws.Range('D37').Value = '0.02326'
$ws.Range('E37').Value = '  -0.72%  '
$ws.Range('D38').Value = '5.181'
$ws.Range('E38').Value = '  -0.97%  '
$ws.Range('D39').Value = '8.750'
$ws.Range('E39').Value = '  -0.29%  '
$ws.Range('D40').Value = '0.6275'
$ws.Range('E40').Value = '  -0.07%  '
$ws.Range('D41').Value = '11.20'
$ws.Range('E41').Value = '  -1.10%  '
$ws.Range('D42').Value = '1.172'
$ws.Range('E42').Value = '  -1.20%  '
$ws.Range('E43').Value = '  -0.06%  '
$ws.Range('D44').Value = '1.392'
$ws.Range('E44').Value = '  -0.54%  '
$ws.Range('D45').Value = '13.50'
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('E46').Value = '  +0.61%  '
$ws.Range('D47').Value = '0.5887'
$ws.Range('E47').Value = '  +0.61%  '
$ws.Range('D48').Value = '124.31'
$ws.Range('E48').Value = '  -1.53%  '
$ws.Range('D49').Value = '1.987'
$ws.Range('E49').Value = '  -0.14%  '
$ws.Range('E50').Value = '  -0.25%  '
$ws.Range('D51').Value = '0.06904'
$ws.Range('E51').Value = '  +0.05%  '
